$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update HZ_Detection_Limit (column Q) values to reflect SOLAR_GRANULATION_RMS = 0.8
$newQ = @{
  2 = 0.8239175227358072
  3 = 1.086836000081502
  4 = 1.673011684280936
  5 = 1.155518363889651
  6 = 0.8995655422072516
  7 = 0.8995655422072516
  8 = 1.639734192087474
  9 = 1.855414867436865
  10 = 1.44974443522691
  11 = 2.511247632757248
  12 = 1.607880801785313
  13 = 1.04713564664033
  14 = 1.326232586610272
  15 = 1.98896214143939
  16 = 0.7360760259646952
  17 = 0.7360760259646952
  18 = 2.674574616230561
  19 = 1.711840798525595
  20 = 0.8422036146915931
  21 = 1.567080541948636
  22 = 0.6036502700178147
  23 = 0.6036502700178147
  24 = 0.6036502700178147
  25 = 1.645503504737809
  26 = 0.8663687244288396
  27 = 1.429141886451972
  28 = 2.167290027842529
  29 = 0.8066475768186987
  30 = 1.390642899182112
  31 = 1.390642899182112
  32 = 1.390642899182112
  33 = 1.390642899182112
  34 = 1.390642899182112
  35 = 0.9117968498801516
  36 = 1.559962907000137
  37 = 1.157049779910267
  38 = 1.157049779910267
  39 = 0.9201520384433114
  40 = 1.516378108953446
  41 = 1.506642040307112
  42 = 0.8164897989250212
  43 = 0.8226369723154777
  44 = 1.535217968277671
  45 = 1.380080961817654
  46 = 0.8627064625966099
  47 = 1.420266750973085
  48 = 1.420266750973085
  49 = 1.234228002760053
  50 = 0.9014348176269065
  51 = 1.417541627563252
  52 = 1.484563724091384
  53 = 0.8839610806037723
  54 = 1.102568633308624
  55 = 1.102568633308624
  56 = 1.033524097456652
  57 = 1.342540001201114
  58 = 1.342540001201114
  59 = 1.626298417461924
  60 = 1.626298417461924
  61 = 1.183025290588741
  62 = 1.183025290588741
  63 = 1.183025290588741
  64 = 1.77076379929166
  65 = 1.822544105549942
  66 = 1.776874771529581
  67 = 1.435841770397525
  68 = 1.589191432361864
  69 = 1.589191432361864
  70 = 2.099615239697236
  71 = 2.012620944277857
  72 = 1.709403743219913
  73 = 1.709403743219913
  74 = 2.345782205913604
  75 = 5.691327750389669
  76 = 5.691327750389669
  77 = 5.691327750389669
  78 = 10.54475952976137
  79 = 10.54475952976137
  80 = 61.86865244577595
}

foreach ($row in $newQ.Keys) {
    $ws.Cells.Item([int]$row, 17).Value = $newQ[$row]
}
